$d = $word.ActiveDocument

# The bullet "Utilized SQL queries ... identifying null values in data
# extensions." gains " and duplicates" right after "null values", becoming
# "... identifying null values and duplicates in data extensions."
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "identifying null values in data extensions.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "identifying null values and duplicates in data extensions.",
    2
)
